# Weekly price-sheet refresh: a new week of "Poroto verde" price quotes is
# prepended to the data block (which starts at row 491), pushing the
# existing rows down by two rows (old row 515 -> new row 517).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data (current rows 491:515) down two rows, making room
# for the two new weekly entries.
$ws.Rows("491:492").Insert()

# New row 491
$ws.Cells.Item(491, 1).Value = 9
$ws.Cells.Item(491, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(491, 3).Value = "Metropolitana"
$ws.Cells.Item(491, 4).Value = 44706
$ws.Cells.Item(491, 5).Value = 13
$ws.Cells.Item(491, 6).Value = 100112031
$ws.Cells.Item(491, 7).Value = "Poroto verde"
$ws.Cells.Item(491, 8).Value = "Magnum"
$ws.Cells.Item(491, 9).Value = "Primera"
$ws.Cells.Item(491, 10).Value = 48
$ws.Cells.Item(491, 11).Value = 35000
$ws.Cells.Item(491, 12).Value = 35000
$ws.Cells.Item(491, 13).Value = 35000
$ws.Cells.Item(491, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(491, 15).Value = "Región Metropolitana"
$ws.Cells.Item(491, 16).Value = 1400
$ws.Cells.Item(491, 17).Value = 25
$ws.Cells.Item(491, 18).Value = "Hortaliza"

# New row 492
$ws.Cells.Item(492, 1).Value = 9
$ws.Cells.Item(492, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(492, 3).Value = "Metropolitana"
$ws.Cells.Item(492, 4).Value = 44706
$ws.Cells.Item(492, 5).Value = 13
$ws.Cells.Item(492, 6).Value = 100112031
$ws.Cells.Item(492, 7).Value = "Poroto verde"
$ws.Cells.Item(492, 8).Value = "Sin especificar"
$ws.Cells.Item(492, 9).Value = "Primera"
$ws.Cells.Item(492, 10).Value = 35
$ws.Cells.Item(492, 11).Value = 45000
$ws.Cells.Item(492, 12).Value = 45000
$ws.Cells.Item(492, 13).Value = 45000
$ws.Cells.Item(492, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(492, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(492, 16).Value = 1800
$ws.Cells.Item(492, 17).Value = 25
$ws.Cells.Item(492, 18).Value = "Hortaliza"
